$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 needs to hold the literal TEXT "1" (not the number 1), matching
# the style it already has (s="23"). Assigning a numeric-looking string
# directly (Value = "1") gets auto-coerced to a numeric cell, and using a
# leading apostrophe (quote-prefixed text) stamps the cell with a brand new
# "quotePrefix" style instead of reusing the existing one. Routing the text
# through a helper cell's TEXT() formula result and pasting only the value
# keeps the original style intact while landing a genuine text cell.
$helper = $ws.Range("Z1")
$helper.Formula = "=TEXT(1,""0"")"
$helper.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$helper.Clear()

